$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.538.66'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '3.496.75'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.61'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.61'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '3.494.02'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  +1.29%  '
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.59'
$ws.Range("E11").Value = '  +6.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '4.086.86'
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '67.485.14'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.486.84'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.43'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.81'
$ws.Range("E21").Value = '  +3.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '446.56'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.628'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.32'
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").Value = '3.639.95'
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000127'
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.72'
$ws.Range("E28").Value = '  +2.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.01'
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("E31").Value = '  +4.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.172'
$ws.Range("E32").Value = '  +4.54%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.57'
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.14'
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("D37").Value = '3.486.94'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.00'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  +7.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '178.53'
$ws.Range("E41").Value = '  +1.03%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0895'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.892'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '30.21'
$ws.Range("E46").Value = '  +7.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.47'
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.31'
$ws.Range("E48").Value = '  +4.81%  '
$ws.Range("E49").Value = '  -3.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.61'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.989'
$ws.Range("E51").Value = '  -0.49%  '
